# Skill Filtering on Project Logic, Data Changes
#
# The "skills" column (G) used to hold a comma-separated list of skill
# names (e.g. "Git,Lambda,Python"). It now holds a JSON array of skill
# UUIDs instead, one array per project row. The "description" column (D)
# is untouched. A stray empty row (row 11, which only ever held a blank
# C11 cell) is also removed from the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$skillsRow2 = @"
[
"a3743fe8-7399-45d2-bbce-ff0fdcaa1df8",
"cc160541-0d68-4db5-9b56-a7a5001b02e3",
"abbecdfe-8116-4ffc-a424-c1b237deb205",
"5d4de0fc-2124-4356-8379-0fcc57c4f247",
"56250a05-8bbd-425b-9704-2282595da788",
"a771ac80-ac0e-4d82-b9ec-4e31f5bacdba",
"bb21bd4c-0280-4eaa-aaf1-ad88893836f5",
"dcdd7b94-0137-456a-b67b-00419e2d4a84",
"6e5df643-1ee5-405a-9021-922da377d330"
]
"@

$skillsRow3 = @"
[
"4424a537-3628-4a37-abb5-d23f2cf45fc2",
"0fed5a25-c2fc-4901-b0fb-eb2b6140b4ac",
"bb21bd4c-0280-4eaa-aaf1-ad88893836f5",
"3e9c896d-a80f-42b7-be71-fe1a3e0f9f0f",
"5dc4d99b-2be7-4413-a515-5403c97f28f2",
"c9b582ed-2a96-412d-99ac-d3cf2283cc5d",
"5d4de0fc-2124-4356-8379-0fcc57c4f247",
"457a4d4d-00e9-4522-8fe7-58044da7d0d7",
"12cd7cf3-5963-4edd-a367-5cfb120e3e55",
"dcdd7b94-0137-456a-b67b-00419e2d4a84",
"c2fdf2be-2dfd-4e45-a4c5-f81a7163dad2",
"cb0a561b-8483-4089-810b-637e4d689ee4",
"b06bb60c-31e1-4939-8d18-54bc6f178cb6",
"56250a05-8bbd-425b-9704-2282595da788"
]
"@

$skillsRow4 = @"
[
"4424a537-3628-4a37-abb5-d23f2cf45fc2",
"10f3de08-5f75-4cfa-b481-05daaee8db79",
"9e9c0bac-7b6b-41fb-acc5-acf3f651c093",
"ed450cf6-a2e6-4512-915e-c112eff7c797",
"5d4de0fc-2124-4356-8379-0fcc57c4f247",
"ac160f48-1c9f-4e5e-8796-8fecbcd04367",
"457a4d4d-00e9-4522-8fe7-58044da7d0d7",
"dcdd7b94-0137-456a-b67b-00419e2d4a84",
"f2444e66-730a-4c4a-b59e-d342e8d57ed6",
"bb21bd4c-0280-4eaa-aaf1-ad88893836f5",
"c2fdf2be-2dfd-4e45-a4c5-f81a7163dad2",
"cb0a561b-8483-4089-810b-637e4d689ee4",
"fa917912-9e24-4330-b9e6-e6b50d5ed0e2"
]
"@

$skillsRow5 = @"
[
"e8a6997b-c4d1-4ffb-a214-4ed99487b858",
"d0eb1478-f755-49be-97bd-1b98c27e4b4a",
"0535d6cc-3815-46e7-bf33-406bb9e3b094",
"115c0855-4c6e-4943-963a-7ff60139a932",
"f56d0a66-2b46-4697-b53e-4aad36ea0fc5",
"fc99bf5d-ead9-44eb-8de7-dd0add714411",
"5d12e67a-f6db-48fd-848c-e267add62128",
"f6a97e57-0b4e-459a-b8d7-0dffe189be5a",
"79de5ad4-24ca-47c1-b18d-798c68705e1b",
"48ec5901-f927-4b3d-8a12-4c0d67185367",
"823b6fa1-5bcf-498e-bf5e-29387b6fac3b",
"6a7cbf58-8987-44ab-8065-fc527ed632e2",
"2017a078-6d1f-4705-af3a-f60e2289a8ee",
"e2a283c2-85f3-4941-b38e-71cedb3e4dec"
]
"@

$ws.Range("G2").Value = $skillsRow2
$ws.Range("G3").Value = $skillsRow3
$ws.Range("G4").Value = $skillsRow4
$ws.Range("G5").Value = $skillsRow5

# Writing a multi-line value makes the engine auto-grow the row height
# (the way Excel would with wrap-text autofit); the source rows used a
# fixed explicit height of 15, so put that back.
$ws.Rows(2).RowHeight = 15
$ws.Rows(3).RowHeight = 15
$ws.Rows(4).RowHeight = 15
$ws.Rows(5).RowHeight = 15

# Remove the stray empty row 11 (only ever held a blank C11 cell) without
# shifting the later, still-sparse row (originally row 18) out of place:
# delete row 11 (shifting everything below it up by one row), then
# re-insert a blank row just above where row 18 now sits so it lands
# back on row 18, unchanged.
$ws.Rows(11).Delete()
$ws.Rows(17).Insert()
